$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -12.406
$ws.Range("B3").Value = 6.247999999999999
$ws.Range("C5").Value = -12.799
$ws.Range("D5").Value = -8.309000000000001
$ws.Range("E7").Value = 13.078
$ws.Range("D9").Value = -7.902000000000001
$ws.Range("D11").Value = -8.121
$ws.Range("E11").Value = 12.775
$ws.Range("B14").Value = 6.331
$ws.Range("B16").Value = 5.914
$ws.Range("C16").Value = -11.856
$ws.Range("D17").Value = -7.828
$ws.Range("E19").Value = 12.756
$ws.Range("B21").Value = 6.434
$ws.Range("D21").Value = -7.890000000000001
$ws.Range("E21").Value = 12.054
$ws.Range("B23").Value = 6.547
$ws.Range("B25").Value = 6.103999999999999
